$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert 4 new rows at the top of the data block (rows 936-939),
# pushing the existing rows 936-948 down to 940-952.
$ws.Rows("936:939").Insert()

# --- Row 936 (new) ---
$ws.Range("A936").Value = 8
$ws.Range("B936").Value = "Terminal La Palmera de La Serena"
$ws.Range("C936").Value = "Coquimbo"
$ws.Range("D936").Value = 44595
$ws.Range("E936").Value = 4
$ws.Range("F936").Value = 100112020
$ws.Range("G936").Value = "Tomate"
$ws.Range("H936").Value = "Semiduro"
$ws.Range("I936").Value = "Primera"
$ws.Range("J936").Value = 800
$ws.Range("K936").Value = 9000
$ws.Range("L936").Value = 10000
$ws.Range("M936").Value = 9500
$ws.Range("N936").Value = "$/bandeja 18 kilos"
$ws.Range("O936").Value = "Provincia de Limarí"
$ws.Range("P936").Value = 528
$ws.Range("Q936").Value = 18
$ws.Range("R936").Value = "Hortaliza"

# --- Row 937 (new) ---
$ws.Range("A937").Value = 8
$ws.Range("B937").Value = "Terminal La Palmera de La Serena"
$ws.Range("C937").Value = "Coquimbo"
$ws.Range("D937").Value = 44595
$ws.Range("E937").Value = 4
$ws.Range("F937").Value = 100112020
$ws.Range("G937").Value = "Tomate"
$ws.Range("H937").Value = "Semiduro"
$ws.Range("I937").Value = "Primera"
$ws.Range("J937").Value = 800
$ws.Range("K937").Value = 5500
$ws.Range("L937").Value = 6000
$ws.Range("M937").Value = 5750
$ws.Range("N937").Value = "$/caja 10 kilos"
$ws.Range("O937").Value = "Provincia del Elquí"
$ws.Range("P937").Value = 575
$ws.Range("Q937").Value = 10
$ws.Range("R937").Value = "Hortaliza"

# --- Row 938 (new) ---
$ws.Range("A938").Value = 8
$ws.Range("B938").Value = "Terminal La Palmera de La Serena"
$ws.Range("C938").Value = "Coquimbo"
$ws.Range("D938").Value = 44595
$ws.Range("E938").Value = 4
$ws.Range("F938").Value = 100112020
$ws.Range("G938").Value = "Tomate"
$ws.Range("H938").Value = "Semiduro"
$ws.Range("I938").Value = "Segunda"
$ws.Range("J938").Value = 400
$ws.Range("K938").Value = 7000
$ws.Range("L938").Value = 7500
$ws.Range("M938").Value = 7250
$ws.Range("N938").Value = "$/bandeja 18 kilos"
$ws.Range("O938").Value = "Provincia de Limarí"
$ws.Range("P938").Value = 403
$ws.Range("Q938").Value = 18
$ws.Range("R938").Value = "Hortaliza"

# --- Row 939 (new) ---
$ws.Range("A939").Value = 8
$ws.Range("B939").Value = "Terminal La Palmera de La Serena"
$ws.Range("C939").Value = "Coquimbo"
$ws.Range("D939").Value = 44595
$ws.Range("E939").Value = 4
$ws.Range("F939").Value = 100112020
$ws.Range("G939").Value = "Tomate"
$ws.Range("H939").Value = "Semiduro"
$ws.Range("I939").Value = "Segunda"
$ws.Range("J939").Value = 400
$ws.Range("K939").Value = 4000
$ws.Range("L939").Value = 4500
$ws.Range("M939").Value = 4250
$ws.Range("N939").Value = "$/caja 10 kilos"
$ws.Range("O939").Value = "Provincia del Elquí"
$ws.Range("P939").Value = 425
$ws.Range("Q939").Value = 10
$ws.Range("R939").Value = "Hortaliza"
